$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("B2").Value = 7190
$ws.Range("C3").Value = 169826
$ws.Range("C4").Value = 160654
$ws.Range("C5").Value = 9172
$ws.Range("C8").Value = 65.63
